# Update the "build_version" / version string throughout the workbook.
#
# Old version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
# New version string: "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$oldVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"
$newVersion = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on February 03 2026 17.29.55 EST)"

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- About sheet ---

# A2: "Version: ..."
$a2 = $wsAbout.Range("A2")
$a2.Value = $a2.Value().Replace($oldVersion, $newVersion)

# A6: Recommended citation text containing the version string
$a6 = $wsAbout.Range("A6")
$a6.Value = $a6.Value().Replace($oldVersion, $newVersion)

# --- Boundaries and methane sources sheet ---
# Column S ("build_version") rows 2 through 8 contain the plain version string
for ($row = 2; $row -le 8; $row++) {
    $cell = $wsData.Cells.Item($row, 19)  # column S = 19
    $cur = $cell.Value()
    if ($cur -eq $oldVersion) {
        $cell.Value = $newVersion
    }
}
